$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values in column D are plain numbers (e.g. 215.27) which Excel
# would otherwise auto-convert to a numeric cell. The source data keeps them as
# plain text (same as the other already-text price cells), so force text format
# on just those cells before writing the new value.
$textCells = @("D5", "D9", "D10", "D11", "D15", "D17", "D21", "D22", "D26", "D29", "D35", "D38", "D42", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.759.68"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.634.68"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "215.27"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "0.0640"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "19.87"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "0.0779"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.632.18"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.860.64"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "0.555"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "0.0₃0776"
$ws.Range("D17").Value = "63.12"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "25.784.37"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "194.04"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "9.95"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D26").Value = "140.65"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").Value = "15.57"
$ws.Range("E29").Value = "  +0.60%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "2.39"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "0.553"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "1.106.00"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("D42").Value = "5.56"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "99.20"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").Value = "0.0₆0108"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("D46").Value = "55.17"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  +11.64%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.418"
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.68"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  -0.53%  "
